$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: label for the new "worst case" block
$ws.Range("C10").Value = "worst case"

# Row 11: raw
$ws.Range("B11").Value = "raw"
$ws.Range("C11").Value = 128000
$ws.Range("C11").NumberFormat = 'General\ "Bytes"'

# Row 12: compressed
$ws.Range("B12").Value = "compressed"
$ws.Range("C12").Value = 114179
$ws.Range("C12").NumberFormat = 'General\ "Bytes"'

# Row 13: ratio
$ws.Range("B13").Value = "ratio"
$ws.Range("C13").Formula = "=(C12/C11)"
$ws.Range("C13").NumberFormat = $ws.Range("F4").NumberFormat

$ws.Range("D16").Select()
